# Fruta / hortaliza, semanal
# Insert the newest weekly price observation for
# "Terminal La Palmera de La Serena - Jengibre" as a new row 153,
# pushing the existing historical rows (153..189) down by one
# (they keep their original data, just shifted), and growing the
# sheet's used range from A1:R189 to A1:R190.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 153..189 down to 154..190, leaving a blank row 153.
$ws.Rows.Item(153).Insert()

# Populate the new row 153 with this week's observation.
$ws.Cells.Item(153, 1).Value = 8
$ws.Cells.Item(153, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(153, 3).Value = "Coquimbo"
$ws.Cells.Item(153, 4).Value = 45211
$ws.Cells.Item(153, 5).Value = 4
$ws.Cells.Item(153, 6).Value = 100114007
$ws.Cells.Item(153, 7).Value = "Jengibre"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 340
$ws.Cells.Item(153, 11).Value = 21000
$ws.Cells.Item(153, 12).Value = 22000
$ws.Cells.Item(153, 13).Value = 21500
$ws.Cells.Item(153, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(153, 15).Value = "Perú"
$ws.Cells.Item(153, 16).Value = 1654
$ws.Cells.Item(153, 17).Value = 13
$ws.Cells.Item(153, 18).Value = "Hortaliza"
